# Mark column "N" (checked) as 1 for the rows that correspond to the
# "mega delivery distance" functionality rollout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 46, 50)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 14).Value = 1
}
